# Edit: remove outlier provinces (Iloilo, Sulu, Tawi-Tawi) from the coastal
# exposure dataset and recompute resilience/risk tiers for the remaining rows.
#
# The underlying workbook lists Philippine provinces (column A) together with
# a "resilience" tier (column B) and a "risk" tier (column C). Three outlier
# provinces are dropped entirely (their rows removed), which shifts every
# subsequent province up by one row; the resilience/risk tiers for the
# remaining 68 provinces are also refreshed to their recalculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the 3 trailing rows so the sheet shrinks from 72 data rows to 69
# (1 header + 68 province rows), matching the new row count after removing
# the outlier provinces.
$ws.Rows("70:72").Delete()

$data = @'
2|Abra|Mid|Mid
3|Agusan Del Norte|Mid|High
4|Agusan Del Sur|Low|High
5|Aklan|Low|High
6|Albay|Mid|Low
7|Antique|Low|High
8|Apayao|Mid|Low
9|Aurora|High|Mid
10|Basilan|Low|Mid
11|Bataan|High|Low
12|Batangas|High|Low
13|Biliran|High|High
14|Bohol|Mid|High
15|Bukidnon|Low|Low
16|Bulacan|High|Low
17|Cagayan|High|Mid
18|Camarines Norte|Mid|High
19|Camarines Sur|Mid|Mid
20|Camiguin|Low|High
21|Capiz|Mid|Mid
22|Catanduanes|High|Mid
23|Cavite|High|Mid
24|Cebu|High|Mid
25|Compostela Valley|Low|Low
26|Davao Del Norte|Mid|Low
27|Davao Del Sur|High|Low
28|Davao Oriental|Low|Mid
29|Eastern Samar|Low|High
30|Guimaras|Mid|High
31|Ifugao|Mid|Low
32|Ilocos Norte|High|Low
33|Ilocos Sur|High|Low
34|Isabela|Mid|Mid
35|Kalinga|High|Low
36|La Union|High|Low
37|Laguna|High|Low
38|Lanao Del Norte|Mid|Mid
39|Lanao Del Sur|Low|Mid
40|Leyte|Mid|High
41|Maguindanao|Low|High
42|Marinduque|Mid|High
43|Masbate|Low|High
44|Misamis Oriental|High|High
45|Negros Occidental|Mid|Mid
46|North Cotabato|Low|Mid
47|Northern Samar|Low|High
48|Nueva Ecija|Mid|Mid
49|Nueva Vizcaya|High|Mid
50|Occidental Mindoro|Mid|High
51|Oriental Mindoro|Mid|High
52|Palawan|Mid|High
53|Pampanga|High|Low
54|Pangasinan|High|Mid
55|Quezon|High|Low
56|Rizal|Low|High
57|Romblon|Low|Mid
58|Samar|Low|Low
59|Sarangani|Low|High
60|Sorsogon|Mid|Low
61|South Cotabato|Low|Mid
62|Southern Leyte|Low|High
63|Sultan Kudarat|Low|High
64|Surigao Del Norte|High|Mid
65|Tarlac|High|Low
66|Zambales|Low|High
67|Zamboanga Del Norte|Mid|Mid
68|Zamboanga Del Sur|Low|Mid
69|Zamboanga Sibugay|Low|Mid
'@

$dataLines = $data -split "`n"

foreach ($line in $dataLines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split '\|'
    $rowNum = [int]$parts[0]
    $province = $parts[1]
    $resilience = $parts[2]
    $risk = $parts[3]

    $ws.Cells.Item($rowNum, 1).Value2 = $province
    $ws.Cells.Item($rowNum, 2).Value2 = $resilience
    $ws.Cells.Item($rowNum, 3).Value2 = $risk
}

Write-Host "Applied outlier removal and tier refresh to $($dataLines.Count) province rows."
